$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.842.97"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.289.24"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'314.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'103.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'39.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "'0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'8.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").Value = "'0.985"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "'15.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "2.633.48"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "2.319.51"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "42.723.88"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.12%  "
$ws.Range("D22").Value = "'73.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "'3.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").Value = "'264.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "'10.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'7.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +18.23%  "
$ws.Range("D30").Value = "'22.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "'36.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'167.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").Value = "'0.0870"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").Value = "'2.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  -5.11%  "
$ws.Range("D37").Value = "'4.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "'0.0350"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").Value = "'3.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("D41").Value = "'1.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("D42").Value = "'70.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "'0.230"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "'94.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.736.10"
$ws.Range("E46").Value = "  +9.38%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "'12.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'80.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'112.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "'8.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.10%  "
